$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.169.54'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '1.794.37'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.16'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5203'
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  -3.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07977'
$ws.Range("E9").Value = '  -2.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.46'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.096'
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.292'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.53'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.293'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").Value = '1.790.90'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.88'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001090'
$ws.Range("E18").Value = '  -3.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06564'
$ws.Range("E19").Value = '  -1.41%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.33'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.955'
$ws.Range("D23").Value = '28.206.93'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.269'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.56'
$ws.Range("E26").Value = '  +3.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.45'
$ws.Range("E27").Value = '  -3.90%  '
$ws.Range("D28").Value = '1.996.62'
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.343'
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.81'
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1075'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.053'
$ws.Range("E32").Value = '  -5.62%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.547'
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07237'
$ws.Range("E35").Value = '  +2.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.07'
$ws.Range("E36").Value = '  +6.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02309'
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2144'
$ws.Range("E38").Value = '  -3.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.690'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.070'
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6167'
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.166'
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.353'
$ws.Range("E43").Value = '  -3.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.21'
$ws.Range("E44").Value = '  -1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.764'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5958'
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.74'
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.223'
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.84'
$ws.Range("E51").Value = '  -1.98%  '
